$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4005110349344403
$ws.Range("C2").Value = 2.655140977822331
$ws.Range("D2").Value = 21.63180862528367
$ws.Range("E2").Value = 4.651000819746614
$ws.Range("F2").Value = 4.689892558235588
$ws.Range("G2").Value = 42

$ws.Range("B3").Value = 0.5668892295343088
$ws.Range("C3").Value = 2.479136195919392
$ws.Range("D3").Value = 21.0318797357638
$ws.Range("E3").Value = 4.586052740185595
$ws.Range("F3").Value = 4.60741567970938
$ws.Range("G3").Value = 41

$ws.Range("B4").Value = 0.3205000134895055
$ws.Range("C4").Value = 2.637395790930121
$ws.Range("D4").Value = 22.34565928591471
$ws.Range("E4").Value = 4.727119554857346
$ws.Range("F4").Value = 4.776323983692677
$ws.Range("G4").Value = 40

$ws.Range("B5").Value = 0.5762346695970005
$ws.Range("C5").Value = 2.589095673025974
$ws.Range("D5").Value = 22.01892221720106
$ws.Range("E5").Value = 4.692432441410431
$ws.Range("F5").Value = 4.717794302558035
$ws.Range("G5").Value = 39

$ws.Range("B6").Value = 0.2753148666760855
$ws.Range("C6").Value = 2.600068871628832
$ws.Range("D6").Value = 23.53914279779983
$ws.Range("E6").Value = 4.851715448972645
$ws.Range("F6").Value = 4.908919327970988
$ws.Range("G6").Value = 38

$ws.Range("B7").Value = 0.5933497739086353
$ws.Range("C7").Value = 2.532836127095716
$ws.Range("D7").Value = 22.7469172165638
$ws.Range("E7").Value = 4.769372832623154
$ws.Range("F7").Value = 4.797596535730608
$ws.Range("G7").Value = 37

$ws.Range("B8").Value = 0.08217900475735471
$ws.Range("C8").Value = 2.189510016720857
$ws.Range("D8").Value = 19.90140921937524
$ws.Range("E8").Value = 4.461099552730833
$ws.Range("F8").Value = 4.52361299942514
$ws.Range("G8").Value = 36

$ws.Range("B9").Value = 0.3064881723285103
$ws.Range("C9").Value = 2.115380499113829
$ws.Range("D9").Value = 19.45264066361694
$ws.Range("E9").Value = 4.410514784423349
$ws.Range("F9").Value = 4.464087741054711
$ws.Range("G9").Value = 35

$ws.Range("B10").Value = 0.08204119957200627
$ws.Range("C10").Value = 2.253188973636722
$ws.Range("D10").Value = 20.48504743802246
$ws.Range("E10").Value = 4.526041033621155
$ws.Range("F10").Value = 4.593350817267506
$ws.Range("G10").Value = 34

$ws.Range("B11").Value = 0.4071412824886091
$ws.Range("C11").Value = 2.043145068801755
$ws.Range("D11").Value = 20.02480196731898
$ws.Range("E11").Value = 4.474908040096352
$ws.Range("F11").Value = 4.525442837904833
$ws.Range("G11").Value = 33
